$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 191.73
$ws.Range("I15").Value = 191.73
$ws.Range("K15").Value = 575.1899999999999
$ws.Range("M15").Value = -406.1899999999999

$ws.Range("H19").Value = 841.8823
$ws.Range("I19").Value = 1135.2
$ws.Range("J19").Value = 719.6667
$ws.Range("K19").Value = 1135.2
$ws.Range("L19").Value = 719.6667
$ws.Range("M19").Value = -960.2
$ws.Range("N19").Value = -1069.6667

$ws.Range("H74").Value = 8367
$ws.Range("I74").Value = 9208.75
$ws.Range("K74").Value = 9208.75
$ws.Range("M74").Value = -8272.75

$ws.Range("H77").Value = 8367
$ws.Range("I77").Value = 9208.75
$ws.Range("K77").Value = 46043.75
$ws.Range("M77").Value = -41363.75

$ws.Range("H96").Value = 1411.8966
$ws.Range("I96").Value = 511.52942
$ws.Range("J96").Value = 2687.4167
$ws.Range("K96").Value = 1534.58826
$ws.Range("L96").Value = 8062.250100000001
$ws.Range("M96").Value = -161.58826
$ws.Range("N96").Value = -10808.2501

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value = 609.46155
$ws.Range("I107").Value = 644.7826
$ws.Range("J107").Value = 338.66666
$ws.Range("K107").Value = 644.7826
$ws.Range("L107").Value = 338.66666
$ws.Range("M107").Value = 1275.2174
$ws.Range("N107").Value = -4178.66666

$ws.Range("H132").Value = 8343291.5
$ws.Range("I132").Value = 9269435
$ws.Range("K132").Value = 27808305
$ws.Range("M132").Value = -27805775

$ws.Range("H137").Value = 1677.091
$ws.Range("I137").Value = 1383.1111
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 4149.3333
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -1599.3333
$ws.Range("N137").Value = -14100

$ws.Range("H138").Value = 3179.2856
$ws.Range("J138").Value = 3595.8235
$ws.Range("L138").Value = 10787.4705
$ws.Range("N138").Value = -21067.4705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32590.014
$ws.Range("I32").Value = 13041.017
$ws.Range("J32").Value = 96667.28
$ws.Range("K32").Value = 13041.017
$ws.Range("L32").Value = 96667.28
$ws.Range("M32").Value = -12754.017
$ws.Range("N32").Value = -97241.28

$ws.Range("H68").Value = 51695
$ws.Range("J68").Value = 51695
$ws.Range("L68").Value = 51695
$ws.Range("N68").Value = -53317

$ws.Range("H71").Value = 51695
$ws.Range("J71").Value = 51695
$ws.Range("L71").Value = 155085
$ws.Range("N71").Value = -163197

$ws.Range("H97").Value = 85707.664
$ws.Range("I97").Value = 201658
$ws.Range("J97").Value = 2886
$ws.Range("K97").Value = 201658
$ws.Range("L97").Value = 2886
$ws.Range("M97").Value = -201162
$ws.Range("N97").Value = -3878

$ws.Range("H102").Value = 86351.336
$ws.Range("I102").Value = 144897.14
$ws.Range("K102").Value = 144897.14
$ws.Range("M102").Value = -143275.14

$ws.Range("H132").Value = 40486.312
$ws.Range("I132").Value = 52398.668
$ws.Range("J132").Value = 4749.25
$ws.Range("K132").Value = 157196.004
$ws.Range("L132").Value = 14247.75
$ws.Range("M132").Value = -154666.004
$ws.Range("N132").Value = -19307.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 421.2
$ws.Range("I7").Value = 421.2
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 421.2
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -308.2
$ws.Range("N7").ClearContents()

$ws.Range("H80").Value = 1095.091
$ws.Range("I80").Value = 679.9231
$ws.Range("J80").Value = 1364.95
$ws.Range("K80").Value = 679.9231
$ws.Range("L80").Value = 1364.95
$ws.Range("M80").Value = 318.0769
$ws.Range("N80").Value = -3360.95

$ws.Range("H83").Value = 1095.091
$ws.Range("I83").Value = 679.9231
$ws.Range("J83").Value = 1364.95
$ws.Range("K83").Value = 3399.6155
$ws.Range("L83").Value = 6824.75
$ws.Range("M83").Value = 1592.3845
$ws.Range("N83").Value = -16808.75

$ws.Range("H94").Value = 582.7143
$ws.Range("I94").Value = 495.8
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 495.8
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -44.80000000000001
$ws.Range("N94").Value = -1702

$ws.Range("H105").Value = 401596.2
$ws.Range("I105").Value = 501990
$ws.Range("J105").Value = 334667
$ws.Range("K105").Value = 501990
$ws.Range("L105").Value = 334667
$ws.Range("M105").Value = -500243
$ws.Range("N105").Value = -338161

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41.909092
$ws.Range("I12").Value = 12.4
$ws.Range("J12").Value = 50.588234
$ws.Range("K12").Value = 37.2
$ws.Range("L12").Value = 151.764702
$ws.Range("M12").Value = 135.8
$ws.Range("N12").Value = -497.764702

$ws.Range("H14").Value = 595.6
$ws.Range("I14").Value = 595.6
$ws.Range("K14").Value = 1786.8
$ws.Range("M14").Value = -1613.8

$ws.Range("H131").Value = 849.42
$ws.Range("I131").Value = 635.8333
$ws.Range("J131").Value = 878.5454999999999
$ws.Range("K131").Value = 1907.4999
$ws.Range("L131").Value = 2635.6365
$ws.Range("M131").Value = 3132.5001
$ws.Range("N131").Value = -12715.6365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 250001800
$ws.Range("I97").Value = 333335330
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 333335330
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -333334834
$ws.Range("N97").Value = -2192

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3854.9565
$ws.Range("I132").Value = 3526.9
$ws.Range("J132").Value = 4107.3076
$ws.Range("K132").Value = 10580.7
$ws.Range("L132").Value = 12321.9228
$ws.Range("M132").Value = -8050.700000000001
$ws.Range("N132").Value = -17381.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2330.1667
$ws.Range("I100").Value = 2167
$ws.Range("J100").Value = 2493.3333
$ws.Range("K100").Value = 2167
$ws.Range("L100").Value = 2493.3333
$ws.Range("M100").Value = -1626
$ws.Range("N100").Value = -3575.3333

$ws.Range("H101").Value = 25340.375
$ws.Range("J101").Value = 25340.375
$ws.Range("L101").Value = 25340.375
$ws.Range("N101").Value = -31830.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6578.617
$ws.Range("I132").Value = 4976
$ws.Range("K132").Value = 14928
$ws.Range("M132").Value = -12398
